$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: F1 changes from "pt_min" to "eta" (reusing the existing shared
# string), which also drops the now-unused "pt_min" shared string on save.
$ws.Range("F1").Value = "eta"

# F2:F14 data values change from a flat placeholder (25) to real pt_min data.
$f_values = @(0.1, 0.3, 0.5, 0.7, 0.89, 1.0900000000000001, 1.29, 1.49, 1.69, 1.91, 2.15, 2.4, 2.63)
for ($i = 0; $i -lt $f_values.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $f_values[$i]
}

# Selection moves from F18 to the whole of column G (matches sqref G1:G1048576).
$ws.Columns("G").Select() | Out-Null
